$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 22717
$ws.Range("B2").Value = "Renan Jesus"
$ws.Range("C2").Value = "Financeiro"
$ws.Range("D2").Value = "Problemas pessoais"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 45101
$ws.Range("G2").Value = 11791.23

# Row 3
$ws.Range("A3").Value = 38374
$ws.Range("B3").Value = "Stephany Rocha"
$ws.Range("C3").Value = "P&D"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 45095
$ws.Range("G3").Value = 5323.49

# Row 4
$ws.Range("A4").Value = 12578
$ws.Range("B4").Value = "Luiz Felipe Rodrigues"
$ws.Range("C4").Value = "P&D"
$ws.Range("D4").Value = "Problemas pessoais"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 45096
$ws.Range("G4").Value = 9684.48

# Row 5
$ws.Range("A5").Value = 44314
$ws.Range("B5").Value = "João Gabriel Azevedo"
$ws.Range("C5").Value = "Atendimento ao Cliente"
$ws.Range("D5").Value = "Outros"
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 45090
$ws.Range("G5").Value = 3951.19

# Row 6
$ws.Range("A6").Value = 46502
$ws.Range("B6").Value = "Laís Santos"
$ws.Range("C6").Value = "Financeiro"
$ws.Range("D6").Value = "Problemas pessoais"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 45104
$ws.Range("G6").Value = 12300.01

# Row 7
$ws.Range("A7").Value = 99299
$ws.Range("B7").Value = "Joana Monteiro"
$ws.Range("C7").Value = "Jurídico"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 45078
$ws.Range("G7").Value = 4493.45

# Row 8
$ws.Range("A8").Value = 59798
$ws.Range("B8").Value = "Ana Julia Ramos"
$ws.Range("C8").Value = "Jurídico"
$ws.Range("D8").Value = "Problemas pessoais"
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 45087
$ws.Range("G8").Value = 10244.34

# Row 9
$ws.Range("A9").Value = 12211
$ws.Range("B9").Value = "Ana Luiza Rodrigues"
$ws.Range("C9").Value = "Financeiro"
$ws.Range("D9").Value = "Consulta médica"
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 45104
$ws.Range("G9").Value = 9987.31

# Row 10
$ws.Range("A10").Value = 74778
$ws.Range("B10").Value = "Dr. Pedro Rodrigues"
$ws.Range("C10").Value = "Engenharia"
$ws.Range("D10").Value = "Viagem de negócios"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 45085
$ws.Range("G10").Value = 2640.57

# Row 11
$ws.Range("A11").Value = 40575
$ws.Range("B11").Value = "Calebe Sales"
$ws.Range("C11").Value = "Jurídico"
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 45079
$ws.Range("G11").Value = 12024.9
